{"js": "// \"nieuwe probleem- en doelstelling\"\n// Rewrites the \"Probleemstelling\" and \"Doelstelling\" paragraphs, and moves\n// the _GoBack bookmark to the new last-edited location (inside the\n// \"Doelstelling\" paragraph, right after \"koffie er zijn\").\n\nconst body = context.document.body;\nconst doc = context.document;\n\n// The _GoBack bookmark marks the last edit point. Remove the stale one\n// (it currently sits before the \"Visuele representatie\" drawing) before\n// re-inserting it at the new edit location below.\ndoc.deleteBookmark(\"_GoBack\");\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet probleemIdx = -1;\nlet doelIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (probleemIdx === -1 && t.indexOf(\"Op dit moment is er een systeem op de NHL\") !== -1) {\n    probleemIdx = i;\n  }\n  if (doelIdx === -1 && t.indexOf(\"Een prototype maken\") !== -1) {\n    doelIdx = i;\n  }\n}\n\nconst newProbleemstelling =\n  \"Op dit moment is er een systeem op de NHL om het aantal koppen koffie dat je hebt aangeschaft te meten. \" +\n  \"Als je 10 koppen koffie hebt aangeschaft krijg je de elfde gratis. Dat doen de kassamedewerkers doormiddel van een stempelkaart die de klanten bij zich moeten houden. \" +\n  \"Alleen veel mensen vergeten hem of raken deze kwijt. \" +\n  \"Het percentage van de klanten die er gebruik van maakt is daardoor lager dat het kan zijn.\";\n\nconst newDoelstelling =\n  \"Een prototype van de kaartlezer met pasjes maken die het ouderwetse stempelkaart systeem kan vervangen waarbij je met een schoolpasje via RF-ID de klant met 1-stap de schoolpas tegen de kaartlezer kan leggen. \" +\n  \"De kassa moet automatisch kijken hoeveel koppen koffie er zijn aangeschaft en het corresponderende aantal bijschrijven op de schoolpas.\";\n\nif (probleemIdx !== -1) {\n  paragraphs.items[probleemIdx].insertText(newProbleemstelling, \"Replace\");\n}\nif (doelIdx !== -1) {\n  paragraphs.items[doelIdx].insertText(newDoelstelling, \"Replace\");\n}\nawait context.sync();\n\n// Re-insert _GoBack right after \"koffie er zijn\" inside the Doelstelling\n// paragraph (matches where Word leaves it after the last typed edit).\nconst anchor = body.search(\"koffie er zijn\", { matchCase: true });\nanchor.load(\"items\");\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const endRange = anchor.items[0].getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"nieuwe probleem- en doelstelling\"\n# Rewrites the \"Probleemstelling\" and \"Doelstelling\" paragraphs, and moves\n# the _GoBack bookmark to the new last-edited location (inside the\n# \"Doelstelling\" paragraph, right after \"koffie er zijn\").\n\n$d = $word.ActiveDocument\n\n$newProbleemstelling = \"Op dit moment is er een systeem op de NHL om het aantal koppen koffie dat je hebt aangeschaft te meten. Als je 10 koppen koffie hebt aangeschaft krijg je de elfde gratis. Dat doen de kassamedewerkers doormiddel van een stempelkaart die de klanten bij zich moeten houden. Alleen veel mensen vergeten hem of raken deze kwijt. Het percentage van de klanten die er gebruik van maakt is daardoor lager dat het kan zijn.\"\n$newDoelstelling = \"Een prototype van de kaartlezer met pasjes maken die het ouderwetse stempelkaart systeem kan vervangen waarbij je met een schoolpasje via RF-ID de klant met 1-stap de schoolpas tegen de kaartlezer kan leggen. De kassa moet automatisch kijken hoeveel koppen koffie er zijn aangeschaft en het corresponderende aantal bijschrijven op de schoolpas.\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*Op dit moment is er een systeem op de NHL*\") {\n        $p.Range.Text = $newProbleemstelling\n    }\n    elseif ($t -like \"*Een prototype maken*\") {\n        $p.Range.Text = $newDoelstelling\n    }\n}\n\n# Re-insert _GoBack right after \"koffie er zijn\" inside the Doelstelling\n# paragraph (matches where Word leaves it after the last typed edit).\n# Bookmarks.Add with an existing name relocates it, so no explicit delete\n# of the stale bookmark (near \"Visuele representatie\") is required.\n$rng = $d.Content\n$rng.Find.Execute(\"koffie er zijn\") | Out-Null\n$bmRange = $d.Range($rng.End, $rng.End)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
